$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 7: repeat the year headers above the new "Equipo" summary block, boxed
# with a border that is open on the bottom (it sits flush on top of row 8).
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "0 Años"
$ws.Range("C7").Value = "1 Años"
$ws.Range("D7").Value = "2 Años"
$ws.Range("E7").Value = "3 Años"
$ws.Range("B7:E7").Borders.LineStyle = 1
$ws.Range("B7:E7").Borders(9).LineStyle = 0

# ---------------------------------------------------------------------------
# Rows 8-10: "Equipo 1" / "Equipo 2" weighted averages plus a blank spacer
# row, all sharing the same thin, all-around border. Reset to the "Normal"
# style right after the formulas go in so they don't inherit the 0.00
# number format from the cells they reference.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Equipo 1"
$ws.Range("B8").Formula = "= (0.2)*B2 + 0.8*B3"
$ws.Range("C8:E8").Formula = "= (0.2)*C2 + 0.8*C3"
$ws.Range("A8:E8").Style = "Normal"

$ws.Range("A9").Value = "Equipo 2"
$ws.Range("B9").Formula = "= (1/3)*B2+(1/3)*B3+(1/3)*B4"
$ws.Range("C9:E9").Formula = "= (1/3)*C2+(1/3)*C3+(1/3)*C4"
$ws.Range("A9:E9").Style = "Normal"

$ws.Range("A10:E10").Value = ""

$ws.Range("A8:E10").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# Row 11: "total" row = Equipo 1 + Equipo 2, with A11 highlighted.
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "total"
$ws.Range("B11").Formula = "=B8+B9"
$ws.Range("C11:E11").Formula = "=C8+C9"
$ws.Range("A11:E11").Style = "Normal"
$ws.Range("A11:E11").Borders.LineStyle = 1
$ws.Range("A11").Interior.ThemeColor = 5

# ---------------------------------------------------------------------------
# Rows 2-4: the per-year numbers were stored as text; turn them into real
# numeric values and bump the display format from 0.00 to 0.000. Done last
# so the new formulas above don't pick up this number format by reference.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = 89.18
$ws.Range("C2").Value = 82.79
$ws.Range("D2").Value = 88.2
$ws.Range("E2").Value = 79.52

$ws.Range("B3").Value = 547.33
$ws.Range("C3").Value = 614.07
$ws.Range("D3").Value = 589.3
$ws.Range("E3").Value = 568.93

$ws.Range("B4").Value = 2165.25
$ws.Range("C4").Value = 1834.14
$ws.Range("D4").Value = 2246.75
$ws.Range("E4").Value = 2483.45

$ws.Range("B2:E4").NumberFormat = "0.000"

# ---------------------------------------------------------------------------
# Final cursor position, matching the saved selection in the workbook.
# ---------------------------------------------------------------------------
$ws.Range("E11").Select()
